$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for the newly added provinces (M DIQ FNIDQ, MEDIOUNA),
# keeping the list in alphabetical order (they land between LARACHE and MEKNES).
$ws.Range("A31:A32").EntireRow.Insert()

# Refreshed province / nb_grappe table (recomputed from the updated parquet
# export) covering all 59 provinces now that the two new ones have been added.
$data = @(
    @("AGADIR IDA OUTANANE", 112),
    @("ASSA ZAG", 42),
    @("AZILAL", 336),
    @("BENI MELLAL", 224),
    @("BENSLIMANE", 112),
    @("BERKANE", 112),
    @("BERRECHID", 168),
    @("BOULEMANE", 224),
    @("CHEFCHAOUEN", 336),
    @("CHTOUKA AIT BAHA", 168),
    @("DRIOUCH", 224),
    @("ERRACHIDIA", 400),
    @("FAHS ANJARA", 112),
    @("FES", 56),
    @("FIGUIG", 168),
    @("FQUIH BEN SALAH", 168),
    @("GUELMIM", 168),
    @("GUERCIF", 112),
    @("HAJEB", 168),
    @("HOCEIMA", 280),
    @("IFRANE", 112),
    @("INEZGANE AIT MELLOUL", 56),
    @("JADIDA", 224),
    @("JERADA", 168),
    @("KENITRA", 280),
    @("KHEMISSET", 224),
    @("KHENIFRA", 168),
    @("KHOURIBGA", 168),
    @("LARACHE", 168),
    @("M DIQ FNIDQ", 56),
    @("MEDIOUNA", 56),
    @("MEKNES", 168),
    @("MIDELT", 352),
    @("MOHAMMADIA", 56),
    @("MOULAY YACOUB", 112),
    @("NADOR", 112),
    @("NOUACEUR", 16),
    @("OUARZAZATE", 112),
    @("OUEZZANE", 224),
    @("OUJDA ANGAD", 112),
    @("SALE", 56),
    @("SEFROU", 168),
    @("SETTAT", 280),
    @("SIDI BENNOUR", 224),
    @("SIDI IFNI", 112),
    @("SIDI KACEM", 280),
    @("SIDI SLIMANE", 112),
    @("SKHIRATE TEMARA", 112),
    @("TAN TAN", 39),
    @("TANGER ASSILAH", 112),
    @("TAOUNATE", 224),
    @("TAOURIRT", 168),
    @("TAROUDANNT", 336),
    @("TATA", 168),
    @("TAZA", 280),
    @("TETOUAN", 112),
    @("TINGHIR", 400),
    @("TIZNIT", 168),
    @("ZAGORA", 240)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Refresh the filter database defined name to cover the new extent of the table.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Feuil1!_FilterDatabase") {
        $n.RefersTo = "=Feuil1!`$A`$1:`$E`$60"
    }
}

# Restore the zoomed selection/active cell over the nb_grappe column.
[void]$ws.Range("B2:B60").Select()

Write-Output "done"
